$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 5262.81347962382
$ws.Range("D2").Value = 4281.84807490595
$ws.Range("E2").Value = 6243.7788843417
$ws.Range("F2").Value = 54.5239802822583
$ws.Range("G2").Value = 44.2375577066344
$ws.Range("H2").Value = 65.5439877236182

$ws.Range("C3").Value = 6038.24107142857
$ws.Range("D3").Value = 2546.52968249249
$ws.Range("E3").Value = 9529.95246036466
$ws.Range("F3").Value = 37.4098787139257
$ws.Range("G3").Value = 13.7963507038835
$ws.Range("H3").Value = 65.92337672856

$ws.Range("C4").Value = 5396.43518518519
$ws.Range("D4").Value = 2447.53158962581
$ws.Range("E4").Value = 8345.33878074456
$ws.Range("F4").Value = 43.0402354120618
$ws.Range("G4").Value = 20.5138173773102
$ws.Range("H4").Value = 69.7772868871902

$ws.Range("C5").Value = 3625.54285714286
$ws.Range("D5").Value = -4149.22267867447
$ws.Range("E5").Value = 11400.3083929602
$ws.Range("F5").Value = 13.924323540996
$ws.Range("G5").Value = -25.1256291923018
$ws.Range("H5").Value = 73.3403747406064

$ws.Range("C6").Value = 3901.41014799154
$ws.Range("D6").Value = 2117.62912208941
$ws.Range("E6").Value = 5685.19117389367
$ws.Range("F6").Value = 32.5755195055895
$ws.Range("G6").Value = 18.9339618162075
$ws.Range("H6").Value = 47.7817446234417

$ws.Range("C7").Value = 3725.05072463768
$ws.Range("D7").Value = 1326.91865787289
$ws.Range("E7").Value = 6123.18279140247
$ws.Range("F7").Value = 27.4953899171181
$ws.Range("G7").Value = 9.87861037646889
$ws.Range("H7").Value = 47.9366584126286

$ws.Range("C8").Value = 6000.05294117647
$ws.Range("D8").Value = 2802.2477047623
$ws.Range("E8").Value = 9197.85817759064
$ws.Range("F8").Value = 39.3509596377869
$ws.Range("G8").Value = 17.2298193799924
$ws.Range("H8").Value = 65.6463351617711

$ws.Range("C9").Value = 5325.15124153499
$ws.Range("D9").Value = 3724.153396757
$ws.Range("E9").Value = 6926.14908631298
$ws.Range("F9").Value = 52.5102947324947
$ws.Range("G9").Value = 36.5251185089837
$ws.Range("H9").Value = 70.3671108541237

$ws.Range("C10").Value = 6942.6836935167
$ws.Range("D10").Value = 5217.34329081175
$ws.Range("E10").Value = 8668.02409622165
$ws.Range("F10").Value = 59.9030515466241
$ws.Range("G10").Value = 43.1983816090779
$ws.Range("H10").Value = 78.5563887427439

$ws.Range("C11").Value = 9346.73093220339
$ws.Range("D11").Value = 7304.51487473163
$ws.Range("E11").Value = 11388.9469896752
$ws.Range("F11").Value = 81.0533294440124
$ws.Range("G11").Value = 60.4074873828573
$ws.Range("H11").Value = 104.356471369212

$ws.Range("C12").Value = 3878.97530864198
$ws.Range("D12").Value = 2603.83404312504
$ws.Range("E12").Value = 5154.11657415891
$ws.Range("F12").Value = 39.8773166230098
$ws.Range("G12").Value = 27.4580370845045
$ws.Range("H12").Value = 53.5067081935502

$ws.Range("C13").Value = 8798.75420875421
$ws.Range("D13").Value = 6254.65384865191
$ws.Range("E13").Value = 11342.8545688565
$ws.Range("F13").Value = 85.7811916924413
$ws.Range("G13").Value = 58.2166269370133
$ws.Range("H13").Value = 118.148066071488

$ws.Range("C14").Value = 4496.8275862069
$ws.Range("D14").Value = 1555.24807449396
$ws.Range("E14").Value = 7438.40709791983
$ws.Range("F14").Value = 37.2953957290218
$ws.Range("G14").Value = 16.3768753089616
$ws.Range("H14").Value = 61.9739801257333
